# Insert a new "Match ID" column at the front of the sheet.
#
# The source data (San Jose Home Miscellaneous) originally started at
# column A with "Player ID". This change inserts a brand-new column A
# that identifies which match the row's stats belong to, pushing every
# existing column one place to the right (B->C, C->D, ... W->X).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts all existing columns
# (values, styles, merged cells) one column to the right automatically.
$ws.Columns.Item(1).Insert()

# Header (row 2 is the visible header row; row 1 holds the hidden
# grouping labels like "Performance" / "Aerial Duels").
$ws.Cells.Item(2, 1).Value = "Match ID "

# Match ID value for every data row (4 through 19) plus the hidden
# aggregate row 20.
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 4
}

# Row 3 is a hidden spacer row with no value, left blank.

# Match the bold "header" style used elsewhere in the sheet for the new
# column (bold font, no border) on rows 2-19.
$ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(19, 1)).Font.Bold = $true

# Update the active selection to reflect the newly inserted column's data.
$ws.Range("A2:A19").Select() | Out-Null
